# Update natmi TPM output: drop the "ECs" sending-cluster rows (old rows 2-4)
# and refresh the remaining rows' stats for the new TPM numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows where Sending cluster = "ECs" (old rows 2,3,4).
# Deleting the same row index three times shifts everything below up.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# After the deletion, the old FAPs-sending rows (previously 5-7) are now
# rows 2-4, and the old MuSCs-sending rows (previously 8-10) are now rows
# 5-7. Refresh the cells whose values changed with the new TPM figures.

# Row 2: FAPs -> ECs
$ws.Range("I2").Value = 0.9017494976312432
$ws.Range("J2").Value = 0.9017494976312432
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 55.42213566666666
$ws.Range("N2").Value = 166.266407
$ws.Range("O2").Value = 0.9848186220994556
$ws.Range("P2").Value = 0.9848186220994556
$ws.Range("Q2").Value = 612.5135461028768
$ws.Range("R2").Value = 5512.621914925891
$ws.Range("S2").Value = 0.8880596977360773
$ws.Range("T2").Value = 0.8880596977360773

# Row 3: FAPs -> FAPs
$ws.Range("I3").Value = 0.9017494976312432
$ws.Range("J3").Value = 0.9017494976312432
$ws.Range("O3").Value = 0.00238111271461604
$ws.Range("P3").Value = 0.00238111271461604
$ws.Range("S3").Value = 0.00214716719420838
$ws.Range("T3").Value = 0.00214716719420838

# Row 4: FAPs -> MuSCs
$ws.Range("I4").Value = 0.9017494976312432
$ws.Range("J4").Value = 0.9017494976312432
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.720354
$ws.Range("N4").Value = 2.161062
$ws.Range("O4").Value = 0.01280026518592835
$ws.Range("P4").Value = 0.01280026518592835
$ws.Range("Q4").Value = 7.961197772008
$ws.Range("R4").Value = 71.650779948072
$ws.Range("S4").Value = 0.01154263270095758
$ws.Range("T4").Value = 0.01154263270095758

# Row 5: MuSCs -> ECs
$ws.Range("G5").Value = 1.204152
$ws.Range("H5").Value = 3.612456
$ws.Range("I5").Value = 0.09825050236875665
$ws.Range("J5").Value = 0.09825050236875667
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 55.42213566666666
$ws.Range("N5").Value = 166.266407
$ws.Range("O5").Value = 0.9848186220994556
$ws.Range("P5").Value = 0.9848186220994556
$ws.Range("Q5").Value = 66.73667550728798
$ws.Range("R5").Value = 600.6300795655919
$ws.Range("S5").Value = 0.09675892436337823
$ws.Range("T5").Value = 0.09675892436337824

# Row 6: MuSCs -> FAPs
$ws.Range("G6").Value = 1.204152
$ws.Range("H6").Value = 3.612456
$ws.Range("I6").Value = 0.09825050236875665
$ws.Range("J6").Value = 0.09825050236875667
$ws.Range("O6").Value = 0.00238111271461604
$ws.Range("P6").Value = 0.00238111271461604
$ws.Range("Q6").Value = 0.161357170768
$ws.Range("R6").Value = 1.452214536912
$ws.Range("S6").Value = 0.0002339455204076598
$ws.Range("T6").Value = 0.0002339455204076598

# Row 7: MuSCs -> MuSCs
$ws.Range("G7").Value = 1.204152
$ws.Range("H7").Value = 3.612456
$ws.Range("I7").Value = 0.09825050236875665
$ws.Range("J7").Value = 0.09825050236875667
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.720354
$ws.Range("N7").Value = 2.161062
$ws.Range("O7").Value = 0.01280026518592835
$ws.Range("P7").Value = 0.01280026518592835
$ws.Range("Q7").Value = 0.867415709808
$ws.Range("R7").Value = 7.806741388272001
$ws.Range("S7").Value = 0.001257632484970767
$ws.Range("T7").Value = 0.001257632484970767
